# Updated users and fixtures list
# Applies:
#  - fills in previously-blank prediction cells for existing rows 5, 7, 13
#  - appends three new user rows (15, 16, 17) with their predictions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5: fill in previously empty prediction cells ----
$row5 = @{
    16 = "['Croatia', 1, 0]"      # P5
    17 = "['Germany', 2, 0]"      # Q5
    18 = "['Switzerland', 0, 1]"  # R5
    19 = "['Draw', 1, 1]"         # S5
    20 = "['England', 0, 2]"      # T5
    21 = "['Spain', 2, 1]"        # U5
    22 = "['Slovakia', 1, 0]"     # V5
    23 = "['Austria', 0, 1]"      # W5
    25 = "['Draw', 1, 1]"         # Y5
    26 = "['Portugal', 0, 2]"     # Z5
    27 = "['Belgium', 2, 0]"      # AA5
}
foreach ($col in $row5.Keys) {
    $ws.Cells.Item(5, $col).Value = $row5[$col]
}

# ---- Row 7: fill in previously empty prediction cells ----
$row7 = @{
    19 = "['Serbia', 1, 2]"   # S7
    20 = "['England', 1, 3]"  # T7
    21 = "['Spain', 2, 1]"    # U7
    22 = "['Draw', 1, 1]"     # V7
}
foreach ($col in $row7.Keys) {
    $ws.Cells.Item(7, $col).Value = $row7[$col]
}

# ---- Row 13: fill in previously empty prediction cells ----
$row13 = @{
    18 = "['Switzerland', 1, 2]"       # R13
    19 = "['Slovenia', 2, 1]"          # S13
    20 = "['England', 1, 2]"           # T13
    21 = "['Spain', 2, 1]"             # U13
    22 = "['Draw', 1, 1]"              # V13
    23 = "['Poland', 2, 1]"            # W13
    24 = "['France', 1, 2]"            # X13
    25 = "['Czech Republic', 1, 2]"    # Y13
    26 = "['Portugal', 1, 3]"          # Z13
    27 = "['Belgium', 2, 0]"           # AA13
    28 = "['Germany', 1, 3]"           # AB13
    29 = "['Scotland', 2, 1]"          # AC13
    30 = "['Spain', 0, 3]"             # AD13
}
foreach ($col in $row13.Keys) {
    $ws.Cells.Item(13, $col).Value = $row13[$col]
}

# ---- Row 15: new user "Khushboo92" ----
$row15 = @{
    1  = "Khushboo92"
    2  = "Khushboosoni.iitd@gmail.com"
    3  = "c6e28565ef1b436f33b3c5df169e53521f49ed3d395647146a4b90b6147a1b5e"
    14 = "['Georgia', 2, 3]"           # N15
    15 = "['Portugal', 3, 1]"          # O15
    16 = "['Croatia', 1, 0]"           # P15
    17 = "['Germany', 3, 1]"           # Q15
    18 = "['Draw', 1, 1]"              # R15
    19 = "['Slovenia', 2, 0]"          # S15
    20 = "['England', 1, 3]"           # T15
    21 = "['Spain', 2, 0]"             # U15
    22 = "['Ukraine', 1, 2]"           # V15
    23 = "['Austria', 1, 2]"           # W15
    24 = "['France', 1, 2]"            # X15
    25 = "['Draw', 1, 1]"              # Y15
    26 = "['Portugal', 1, 2]"          # Z15
    27 = "['Belgium', 3, 0]"           # AA15
    28 = "['Germany', 1, 3]"           # AB15
    30 = "['Spain', 1, 2]"             # AD15
    31 = "['Croatia', 2, 1]"           # AE15
    32 = "['Netherlands', 3, 1]"       # AF15
    33 = "['France', 2, 1]"            # AG15
    34 = "['England', 2, 0]"           # AH15
    35 = "['Denmark', 1, 0]"           # AI15
    36 = "['Romania', 1, 2]"           # AJ15
    37 = "['Belgium', 1, 3]"           # AK15
    38 = "['Portugal', 1, 4]"          # AL15
    39 = "['Czech Republic', 3, 2]"    # AM15
}
foreach ($col in $row15.Keys) {
    $ws.Cells.Item(15, $col).Value = $row15[$col]
}

# ---- Row 16: new user "pratham" ----
$row16 = @{
    1  = "pratham"
    2  = "prathammehta@outlook.com"
    3  = "34163a452ad96e0a3882b446d36406363d41d122459694518fa42d1a24aac001"
    15 = "['Portugal', 2, 1]"  # O16
    16 = "['Croatia', 2, 0]"   # P16
    17 = "['Germany', 3, 1]"  # Q16
}
foreach ($col in $row16.Keys) {
    $ws.Cells.Item(16, $col).Value = $row16[$col]
}

# ---- Row 17: new user "Eechalna" ----
$row17 = @{
    1  = "Eechalna"
    2  = "Eechalna@gmail.com"
    3  = "2d8b1fed2294a9f3f766829fb2618c619e6ae5be051c3dbf697d858f60825402"
    21 = "['Spain', 2, 1]"  # U17
}
foreach ($col in $row17.Keys) {
    $ws.Cells.Item(17, $col).Value = $row17[$col]
}
